# Weekly refresh: insert 2 new rows at the top of the "Ají" (chili pepper)
# data block (rows 494-495) and push the existing 90 rows (494-583) down to
# (496-585), matching the original Excel auto-shift behaviour on row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 494; everything that was at
# 494-583 automatically moves down to 496-585 (values + styles included).
$ws.Range("494:495").Insert()

# --- New row 494 -----------------------------------------------------
$ws.Range("A494").Value = 8
$ws.Range("B494").Value = "Terminal La Palmera de La Serena"
$ws.Range("C494").Value = "Coquimbo"
$ws.Range("D494").Value = 45218
$ws.Range("E494").Value = 4
$ws.Range("F494").Value = 100112021
$ws.Range("G494").Value = "Ají"
$ws.Range("H494").Value = "Inferno"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 460
$ws.Range("K494").Value = 25000
$ws.Range("L494").Value = 26000
$ws.Range("M494").Value = 25500
$ws.Range("N494").Value = "$/caja 10 kilos"
$ws.Range("O494").Value = "Región de Arica y Parinacota"
$ws.Range("P494").Value = 2550
$ws.Range("Q494").Value = 10
$ws.Range("R494").Value = "Hortaliza"

# --- New row 495 -----------------------------------------------------
$ws.Range("A495").Value = 8
$ws.Range("B495").Value = "Terminal La Palmera de La Serena"
$ws.Range("C495").Value = "Coquimbo"
$ws.Range("D495").Value = 45218
$ws.Range("E495").Value = 4
$ws.Range("F495").Value = 100112021
$ws.Range("G495").Value = "Ají"
$ws.Range("H495").Value = "Inferno"
$ws.Range("I495").Value = "Segunda"
$ws.Range("J495").Value = 360
$ws.Range("K495").Value = 14000
$ws.Range("L495").Value = 15000
$ws.Range("M495").Value = 14500
$ws.Range("N495").Value = "$/caja 10 kilos"
$ws.Range("O495").Value = "Región de Arica y Parinacota"
$ws.Range("P495").Value = 1450
$ws.Range("Q495").Value = 10
$ws.Range("R495").Value = "Hortaliza"
